$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# Change the target user from "aoliveira" to "ronan" for the grant-statement block (rows 22-52)
$ws.Range("B22:B52").Value2 = "ronan"

# Row 22: keep its own (non-shared) concatenation formula, value recalculates automatically
$ws.Range("D22").Formula = '=A22&" "&B22&" "&C22'

# Rows 23-52: unify into a single concatenation formula "A & " " & B & " " & C"
# (row 23 previously lacked the middle space; rows 36-52 previously held static text)
$ws.Range("D23:D52").Formula = '=A23&" "&B23&" "&C23'

# Row 53 stays blank (mirrors the blank row 21 gap above the grant list).
# Rows 54-83: new blank template rows that only carry the concatenation formula
# (A/B/C are empty, so the cached result is two spaces).
$ws.Range("D54:D83").Formula = '=A54&" "&B54&" "&C54'
